$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Cyborg (01h)" column header to "PuxaFacil (01h)"
$ws.Range("B3").Value = "PuxaFacil (01h)"

# Copy the "X (01h)" cell formatting (highlighted fill) from B7 onto B9, then set its value
$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "X (01h)"

# Update the active cell selection to B9
$ws.Range("B9").Select() | Out-Null
